# Rename the three logo pictures (first-page header BTec logo + both
# footer Pearson logos) exactly as the source commit does: the BTec logo
# picture is renumbered from "image2.jpg" to "image1.jpg", and both
# Pearson logo pictures are renumbered from "image1.png" to "image2.png".
#
# These names are the InlineShape's display Name (<wp:docPr name="...">
# / <pic:cNvPr name="...">), not the media part filenames, so the
# embedded picture bytes/relationships are untouched - only the shape's
# Name property changes.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# First-page header holds the BTec_Logo-Orange inline picture.
$hdr = $sec.Headers.Item(2)   # wdHeaderFooterFirstPage
if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -ge 1) {
    $btecLogo = $hdr.Range.InlineShapes.Item(1)
    $btecLogo.Name = "image1.jpg"
}

# Both the primary and first-page footers hold a Pearson logo inline
# picture; rename each one found.
for ($i = 1; $i -le 3; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        for ($j = 1; $j -le $ftr.Range.InlineShapes.Count; $j++) {
            $pearsonLogo = $ftr.Range.InlineShapes.Item($j)
            $pearsonLogo.Name = "image2.png"
        }
    }
}
